$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D2 value to 30 (was blank)
$ws.Range("D2").Value = 30

# Update the selection to reflect the new active cell E7
$ws.Range("E7").Select()
